$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number format on cells whose new values look numeric,
# so Excel stores them as text (preserving e.g. trailing zeros) instead of
# auto-converting to a Double.
$textCells = @("D5","D6","D9","D11","D12","D14","D19","D20","D21","D22","D23","D24","D25","D27","D29","D31","D32","D33","D35","D36","D37","D39","D40","D41","D42","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '60.495.79'
$ws.Range('E2').Value = '  -3.95%  '
$ws.Range('D3').Value = '2.987.36'
$ws.Range('E3').Value = '  -5.61%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '561.73'
$ws.Range('E5').Value = '  -5.02%  '
$ws.Range('D6').Value = '125.67'
$ws.Range('E6').Value = '  -6.63%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = '2.983.30'
$ws.Range('E8').Value = '  -5.61%  '
$ws.Range('D9').Value = '0.495'
$ws.Range('E9').Value = '  -3.73%  '
$ws.Range('E10').Value = '  -5.70%  '
$ws.Range('D11').Value = '5.10'
$ws.Range('E11').Value = '  -2.46%  '
$ws.Range('D12').Value = '0.436'
$ws.Range('E12').Value = '  -3.83%  '
$ws.Range('E13').Value = '  -6.08%  '
$ws.Range('D14').Value = '32.55'
$ws.Range('E14').Value = '  -6.34%  '
$ws.Range('E15').Value = '  +0.17%  '
$ws.Range('D16').Value = '3.470.38'
$ws.Range('E16').Value = '  -5.87%  '
$ws.Range('D17').Value = '60.620.77'
$ws.Range('E17').Value = '  -3.65%  '
$ws.Range('D18').Value = '2.981.05'
$ws.Range('E18').Value = '  -5.86%  '
$ws.Range('D19').Value = '6.13'
$ws.Range('E19').Value = '  -6.58%  '
$ws.Range('D20').Value = '429.69'
$ws.Range('E20').Value = '  -6.68%  '
$ws.Range('D21').Value = '13.02'
$ws.Range('E21').Value = '  -6.13%  '
$ws.Range('D22').Value = '0.658'
$ws.Range('E22').Value = '  -5.91%  '
$ws.Range('D23').Value = '7.11'
$ws.Range('E23').Value = '  -6.56%  '
$ws.Range('D24').Value = '12.89'
$ws.Range('E24').Value = '  -3.57%  '
$ws.Range('D25').Value = '78.44'
$ws.Range('E25').Value = '  -6.22%  '
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').Value = '  -0.32%  '
$ws.Range('E28').Value = '  -7.33%  '
$ws.Range('D29').Value = '7.17'
$ws.Range('E29').Value = '  -6.99%  '
$ws.Range('E30').Value = '  -6.96%  '
$ws.Range('D31').Value = '25.23'
$ws.Range('E31').Value = '  -6.91%  '
$ws.Range('D32').Value = '6.00'
$ws.Range('E32').Value = '  -10.82%  '
$ws.Range('D33').Value = '0.0932'
$ws.Range('E33').Value = '  -9.81%  '
$ws.Range('E34').Value = '  -4.60%  '
$ws.Range('D35').Value = '0.948'
$ws.Range('E35').Value = '  -8.28%  '
$ws.Range('D36').Value = '5.52'
$ws.Range('E36').Value = '  -5.09%  '
$ws.Range('D37').Value = '49.37'
$ws.Range('E37').Value = '  -3.41%  '
$ws.Range('D38').Value = '0.0₃0659'
$ws.Range('E38').Value = '  -6.22%  '
$ws.Range('D39').Value = '0.0358'
$ws.Range('E39').Value = '  -7.87%  '
$ws.Range('D40').Value = '7.77'
$ws.Range('D41').Value = '375.15'
$ws.Range('E41').Value = '  -6.70%  '
$ws.Range('D42').Value = '0.107'
$ws.Range('E42').Value = '  -4.66%  '
$ws.Range('D43').Value = '2.666.92'
$ws.Range('E43').Value = '  -4.69%  '
$ws.Range('E44').Value = '  -7.44%  '
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('E46').Value = '  -6.49%  '
$ws.Range('D47').Value = '119.46'
$ws.Range('E47').Value = '  -3.18%  '
$ws.Range('D48').Value = '1.96'
$ws.Range('E48').Value = '  -7.20%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').Value = '0.106'
$ws.Range('E49').Value = '  -4.96%  '
$ws.Range('B50').Value = 'Arweave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D50').Value = '32.59'
$ws.Range('E50').Value = '  -4.85%  '
$ws.Range('D51').Value = '23.34'
$ws.Range('E51').Value = '  -7.92%  '

# Restore default ("Normal") style on the text-forced cells so no stray
# number-format styling is left behind on the workbook.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
